$d = $word.ActiveDocument
$paras = $d.Paragraphs

# 1. Date change in the first paragraph (⚡️🚀 banner line)
$paras.Item(1).Range.Text = "⚡️🚀המאמר היומי של מייק 19.06.24:⚡️🚀"

# 2. Paper title
$paras.Item(2).Range.Text = "INTRINSIC DIMENSIONALITY EXPLAINS THE EFFECTIVENESS OF LANGUAGE MODEL FINE-TUNING"

# 3. Paragraph 3 full replacement
$paras.Item(3).Range.Text = "כולכם מכירים את LoRa נכון? בטח גם שמעתם על עשרות השכלולים השונים שלה כמו DoRa, MoRa, GaloRe וכדומה. מתברר כי היה מאמר שבצורה מסוימת הניח יסודות של משפחת הגישות הזו."

# 4. Paragraph 4 full replacement
$paras.Item(4).Range.Text = "למעשה מה זה LoRa? זה אופן שבו אנחנו עושים פיינטיון של מודלים מאומנים גדולים למשימה ספציפית בלי לעדכן את כל משקלי המודל. במקרה של LoRa אנו מאמנים מטריצת תוספות למשקלים של כל שכבה כאשר תוספת זו היא בעלת ראנק נמוך הרבה יותר ממטריצת המשקלים המקורית. כלומר ניתן לייצג אותה על ידי מכפלה שתי מטריצות בעלות רנק נמוך (בגדלים מסוימים במקרה של LoRa)."

# 5. Paragraph 5 full replacement
$paras.Item(5).Range.Text = "מתברר שגישה זו היתה ידוע כבר ב 2020 ואפילו היו מאמרים שדיברו עליה ב 2018. אז המאמרים הציעו מספר דרכים לבניית מטריצת תוספת זו וביניהם הטלה ספארסית של וקטור במימד נמוך למרחב בעל מספר מימדים גבוה דרך Fastfood algorithm (צורה של מטריצת ההטלה הזו - תקראו עליו, זה חמוד)."

# 6. Paragraph 6 full replacement
$paras.Item(6).Range.Text = "בקיצור מאמר ""היסטורי"" מעניין וקל לקריאה."

# 7. Paragraph 7 full replacement (becomes the new link)
$paras.Item(7).Range.Text = "https://arxiv.org/abs/2012.13255"

# 8. Remove the now-obsolete trailing paragraph that held the old link
$paras.Item(8).Range.Delete()
